$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.006697333333333
$ws.Range("H2").Value = 3.020092
$ws.Range("I2").Value = 0.0001985651645046208
$ws.Range("J2").Value = 0.0001985651645046208
$ws.Range("M2").Value = 0.428743
$ws.Range("N2").Value = 1.286229
$ws.Range("O2").Value = 0.00412050394863168
$ws.Range("P2").Value = 0.00412050394863168
$ws.Range("Q2").Value = 0.4316144347853334
$ws.Range("R2").Value = 3.884529913068
$ws.Range("S2").Value = 0.0000008181885444019891
$ws.Range("T2").Value = 0.0000008181885444019891
$ws.Range("G3").Value = 1.006697333333333
$ws.Range("H3").Value = 3.020092
$ws.Range("I3").Value = 0.0001985651645046208
$ws.Range("J3").Value = 0.0001985651645046208
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("O3").Value = 0.7710272268990069
$ws.Range("P3").Value = 0.7710272268990069
$ws.Range("Q3").Value = 80.76353885126798
$ws.Range("R3").Value = 726.8718496614119
$ws.Range("S3").Value = 0.0001530991481467429
$ws.Range("T3").Value = 0.0001530991481467429
$ws.Range("G4").Value = 1.006697333333333
$ws.Range("H4").Value = 3.020092
$ws.Range("I4").Value = 0.0001985651645046208
$ws.Range("J4").Value = 0.0001985651645046208
$ws.Range("M4").Value = 23.39612766666667
$ws.Range("N4").Value = 70.188383
$ws.Range("O4").Value = 0.2248522691523614
$ws.Range("P4").Value = 0.2248522691523614
$ws.Range("Q4").Value = 23.55281933235956
$ws.Range("R4").Value = 211.975373991236
$ws.Range("S4").Value = 0.00004464782781347592
$ws.Range("T4").Value = 0.00004464782781347592
$ws.Range("I5").Value = 0.9806494927176636
$ws.Range("J5").Value = 0.9806494927176637
$ws.Range("M5").Value = 0.428743
$ws.Range("N5").Value = 1.286229
$ws.Range("O5").Value = 0.00412050394863168
$ws.Range("P5").Value = 0.00412050394863168
$ws.Range("Q5").Value = 2131.604894432571
$ws.Range("R5").Value = 19184.44404989314
$ws.Range("S5").Value = 0.004040770106966787
$ws.Range("T5").Value = 0.004040770106966787
$ws.Range("I6").Value = 0.9806494927176636
$ws.Range("J6").Value = 0.9806494927176637
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("O6").Value = 0.7710272268990069
$ws.Range("P6").Value = 0.7710272268990069
$ws.Range("Q6").Value = 398865.1463723194
$ws.Range("R6").Value = 3589786.317350875
$ws.Range("S6").Value = 0.7561074589300181
$ws.Range("T6").Value = 0.7561074589300181
$ws.Range("I7").Value = 0.9806494927176636
$ws.Range("J7").Value = 0.9806494927176637
$ws.Range("M7").Value = 23.39612766666667
$ws.Range("N7").Value = 70.188383
$ws.Range("O7").Value = 0.2248522691523614
$ws.Range("P7").Value = 0.2248522691523614
$ws.Range("Q7").Value = 116319.8005449324
$ws.Range("R7").Value = 1046878.204904391
$ws.Range("S7").Value = 0.2205012636806788
$ws.Range("T7").Value = 0.2205012636806788
$ws.Range("G8").Value = 97.097641
$ws.Range("H8").Value = 291.292923
$ws.Range("I8").Value = 0.01915194211783179
$ws.Range("J8").Value = 0.01915194211783179
$ws.Range("M8").Value = 0.428743
$ws.Range("N8").Value = 1.286229
$ws.Range("O8").Value = 0.00412050394863168
$ws.Range("P8").Value = 0.00412050394863168
$ws.Range("Q8").Value = 41.629933895263
$ws.Range("R8").Value = 374.669405057367
$ws.Range("S8").Value = 0.00007891565312049127
$ws.Range("T8").Value = 0.00007891565312049127
$ws.Range("G9").Value = 97.097641
$ws.Range("H9").Value = 291.292923
$ws.Range("I9").Value = 0.01915194211783179
$ws.Range("J9").Value = 0.01915194211783179
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("O9").Value = 0.7710272268990069
$ws.Range("P9").Value = 0.7710272268990069
$ws.Range("Q9").Value = 7789.778359006915
$ws.Range("R9").Value = 70108.00523106224
$ws.Range("S9").Value = 0.01476666882084214
$ws.Range("T9").Value = 0.01476666882084214
$ws.Range("G10").Value = 97.097641
$ws.Range("H10").Value = 291.292923
$ws.Range("I10").Value = 0.01915194211783179
$ws.Range("J10").Value = 0.01915194211783179
$ws.Range("M10").Value = 23.39612766666667
$ws.Range("N10").Value = 70.188383
$ws.Range("O10").Value = 0.2248522691523614
$ws.Range("P10").Value = 0.2248522691523614
$ws.Range("Q10").Value = 2271.708804968168
$ws.Range("R10").Value = 20445.37924471351
$ws.Range("S10").Value = 0.00430635764386916
$ws.Range("T10").Value = 0.00430635764386916
